$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value to a literal text string without letting Excel's
# autodetect turn date-looking text (e.g. "01/01/2021") into a real date
# serial number + date number-format. We put a formula that evaluates to the
# quoted text literal, then immediately collapse it back down to a plain
# value (Copy + PasteSpecial values-only) so the stored cell stays a plain
# text string (shared string) with its original, untouched style/format.
function Set-TextValue($cell, $text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = "=""$escaped"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

# Ativação: 01/01/2018 -> 01/01/2021
Set-TextValue $ws.Range("B8") "01/01/2021"
Set-TextValue $ws.Range("C8") "01/01/2021"

# Docentes responsáveis: 5840535 - Messias Borges Silva -> 198273 - Domingos Savio Giordani
$ws.Range("B13").Value = "198273 - Domingos Savio Giordani"
$ws.Range("C13").Value = "198273 - Domingos Savio Giordani"

# Método:
$metodo = "Desenvolvimento e apresentação do Projeto monografia a ser desenvolvida na disciplina de Trabalho de Graduação em Engenharia de Produção II, conforme norma do Curso de Engenharia de Produção"
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Critério:
$criterio = "Avaliação Ad hoc por 2 examinadores. A nota da disciplina será a média das duas notas"
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Norma de recuperação:
$norma = "Reapresentação do trabalho modificado para nova avaliação"
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# Bibliografia:
$bibliografia = "Cauchick-Miguel, P. A.,   Metodologia de pesquisa em engenharia de produção e gestão de operações / Afonso Fleury ... [et al.] ; coordenação . - 3. ed. - Rio de Janeiro : Elsevier, 2018. Cauchick-Miguel, P. A.,   Metodologia de pesquisa em engenharia , 1. ed. - Rio de Janeiro : GEN LTC, 2019. BOOTH, W.; COLOMB, G.; WILLIAMS, J. A arte da Pesquisa. 3 ed. Martins Fontes. São Paulo. 2005.GIL, A.C. Como elaborar projetos de pesquisa. 5 ed. Atlas, São Paulo, 2010.MEDEIROS, J. B. Redação Cientifica: A Prática de Fichamentos, Resumos e Resenhas. 11 ed. São Paulo: Atlas, 2009"
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
